# Auto update Excel log
# Appends a new alert row (row 16) to the ALERTS sheet, mirroring the
# existing FALL_DETECTED / CRITICAL rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERTS")

$row = 16

# Columns B-F are plain text already (no auto-conversion risk).
$ws.Cells.Item($row, 2).Value = "11:26:40"
$ws.Cells.Item($row, 3).Value = "11:00"
$ws.Cells.Item($row, 4).Value = "Living Room"
$ws.Cells.Item($row, 5).Value = "CRITICAL"
$ws.Cells.Item($row, 6).Value = "FALL_DETECTED"

# Column A holds a date-formatted string ("2026-02-01") that must stay a
# literal text value like the rest of the log (matching the other rows,
# which store it as inline/shared text, not a real date serial). Writing
# it straight to .Value would make Excel auto-convert it into a date. Route
# it through a text formula and paste the computed value back so the cell
# ends up holding plain text without picking up any extra number formatting.
$ws.Cells.Item($row, 1).Formula = "=""2026-02-01"""
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)  # xlPasteValues
